# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Granada" (Vega Modelo de Temuco)
# at row 164, pushing the existing rows 164-194 down to 165-195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 164 (shifts rows 164..194 -> 165..195,
# and extends the sheet dimension from A1:T194 to A1:T195).
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(164, 1).Value  = 10
$ws.Cells.Item(164, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(164, 3).Value  = "La Araucanía"
$ws.Cells.Item(164, 4).Value  = 44889
$ws.Cells.Item(164, 5).Value  = 9
$ws.Cells.Item(164, 6).Value  = "Fruta"
$ws.Cells.Item(164, 7).Value  = 100104
$ws.Cells.Item(164, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(164, 9).Value  = 100104001
$ws.Cells.Item(164, 10).Value = "Granada"
$ws.Cells.Item(164, 11).Value = "Wonderfull"
$ws.Cells.Item(164, 12).Value = "Primera"
$ws.Cells.Item(164, 13).Value = 50
$ws.Cells.Item(164, 14).Value = 18000
$ws.Cells.Item(164, 15).Value = 18000
$ws.Cells.Item(164, 16).Value = 18000
$ws.Cells.Item(164, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(164, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(164, 19).Value = 1200
$ws.Cells.Item(164, 20).Value = 15
